$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header-ish cells added above the title row (row 1)
$ws.Range("C1").Value = "user1"
$ws.Range("E1").Value = "password1"

# Remove the extra "CPSC 4899" entry (and its credit value) that used to
# live at E4/F4 in the Summer 2022 block
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null

# Insert a new course ("CYBR 4125") into the Spring 2023 column, shifting
# the existing Spring 2023 rows down by one
$ws.Range("C5").Value = "CYBR 4125"
$ws.Range("C6").Value = "CPSC 4135"
$ws.Range("C7").Value = "CPSC 4148"
$ws.Range("C8").Value = "CPSC 4155"

# Split the old single "CPSC 4175" row (13) for Fall 2024 into two rows:
# row 13 now holds "CPSC 4157" and a new row 14 holds "CPSC 4175"
$ws.Range("A13").Value = "CPSC 4157"
$ws.Range("A14").Value = "CPSC 4175"
$ws.Range("B14").Value = 3
